# "Generate Report for Handback"
# Rows for 38033bfe-6d7f-4df0-a5c9-e64f359b4f38 and 741cf1f0-380b-4407-aacf-b3967867eddf
# move from "Ready for handoff" to "Handed back: in sync with en-US", and their
# per-locale handback target file / handback xliff / handback datetime get filled in.

$wb = $excel.ActiveWorkbook

$statusReady    = "Ready for handoff"
$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: the zh-cn / de-de status columns for the two rows flip from
# "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $statusHandedBack
$wsOverview.Range("F4").Value = $statusHandedBack
$wsOverview.Range("E5").Value = $statusHandedBack
$wsOverview.Range("F5").Value = $statusHandedBack

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 4: 38033bfe-6d7f-4df0-a5c9-e64f359b4f38
$wsZhCn.Range("C4").Value = $statusHandedBack
$wsZhCn.Range("I4").Value = "38033bfe-6d7f-4df0-a5c9-e64f359b4f38.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/313a6df70a547950e789c36f699af4fa94cc03a7/e2e/38033bfe-6d7f-4df0-a5c9-e64f359b4f38.md", "", "", "38033bfe-6d7f-4df0-a5c9-e64f359b4f38.md") | Out-Null
$wsZhCn.Range("J4").Value = "38033bfe-6d7f-4df0-a5c9-e64f359b4f38.4a68f86fa29f860ca8a2957390828e9bb0a7a9ac.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-17 04:23:28"

# Row 5: 741cf1f0-380b-4407-aacf-b3967867eddf
$wsZhCn.Range("C5").Value = $statusHandedBack
$wsZhCn.Range("I5").Value = "741cf1f0-380b-4407-aacf-b3967867eddf.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/313a6df70a547950e789c36f699af4fa94cc03a7/e2e/741cf1f0-380b-4407-aacf-b3967867eddf.md", "", "", "741cf1f0-380b-4407-aacf-b3967867eddf.md") | Out-Null
$wsZhCn.Range("J5").Value = "741cf1f0-380b-4407-aacf-b3967867eddf.364ab13b00139d54b41264e1eddc2d21daab8e3f.zh-cn.xlf"
$wsZhCn.Range("K5").Value = "2016-08-17 04:23:28"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4: 38033bfe-6d7f-4df0-a5c9-e64f359b4f38
$wsDeDe.Range("C4").Value = $statusHandedBack
$wsDeDe.Range("I4").Value = "38033bfe-6d7f-4df0-a5c9-e64f359b4f38.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/313a6df70a547950e789c36f699af4fa94cc03a7/e2e/38033bfe-6d7f-4df0-a5c9-e64f359b4f38.md", "", "", "38033bfe-6d7f-4df0-a5c9-e64f359b4f38.md") | Out-Null
$wsDeDe.Range("J4").Value = "38033bfe-6d7f-4df0-a5c9-e64f359b4f38.4a68f86fa29f860ca8a2957390828e9bb0a7a9ac.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-17 04:23:35"

# Row 5: 741cf1f0-380b-4407-aacf-b3967867eddf
$wsDeDe.Range("C5").Value = $statusHandedBack
$wsDeDe.Range("I5").Value = "741cf1f0-380b-4407-aacf-b3967867eddf.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/313a6df70a547950e789c36f699af4fa94cc03a7/e2e/741cf1f0-380b-4407-aacf-b3967867eddf.md", "", "", "741cf1f0-380b-4407-aacf-b3967867eddf.md") | Out-Null
$wsDeDe.Range("J5").Value = "741cf1f0-380b-4407-aacf-b3967867eddf.364ab13b00139d54b41264e1eddc2d21daab8e3f.de-de.xlf"
$wsDeDe.Range("K5").Value = "2016-08-17 04:23:35"
